# Add two new machine master records (rows 31 and 32) with new MAC addresses,
# mirroring the existing data pattern in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: Machine 30 ---
$ws.Range("A31").Value2 = 10030
$ws.Range("B31").Value2 = "Machine 30"
$ws.Range("C31").Value2 = "70-5A-0F-8C-01-39"
$ws.Range("D31").Value2 = "FB5962911663"
$ws.Range("E31").Value2 = "192.168.0.356"
$ws.Range("F31").Value2 = 1001
$ws.Range("G31").Value2 = "eng"
$ws.Range("H31").Value2 = $true
$ws.Range("I31").Value2 = "superadmin"
$ws.Range("J31").Value2 = "now()"
$ws.Range("K31").Value2 = "now()"

# --- Row 32: Machine 31 ---
# (Note: cell-write order intentionally matches the original authoring order
# so the shared-string table ends up with the same ordering: ip before mac.)
$ws.Range("A32").Value2 = 10031
$ws.Range("B32").Value2 = "Machine 31"
$ws.Range("D32").Value2 = "FB5962911663"
$ws.Range("E32").Value2 = "192.168.0.357"
$ws.Range("C32").Value2 = "58-20-B1-DA-F3-FB"
$ws.Range("F32").Value2 = 1001
$ws.Range("G32").Value2 = "eng"
$ws.Range("H32").Value2 = $true
$ws.Range("I32").Value2 = "superadmin"
$ws.Range("J32").Value2 = "now()"
$ws.Range("K32").Value2 = "now()"

# Match the resulting selection seen in the saved workbook (columns beyond
# the data, starting at the row after the header).
$ws.Range("L1:XFD1048576").Select() | Out-Null
